# Adjust values to percents
# - Emp sheet: add 4 new rows (C01..C04 variable codes + new rate labels)
# - New PubAssist sheet added after Food, with public-assistance household data
# - Various sheet selections / active sheet updated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Emp sheet ("Emp"): new rows 24-27.
#    Written bottom-up (27 -> 24) so the shared-string table gets the new
#    strings appended in the same order the reference workbook uses.
# ---------------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item("Emp")

$wsEmp.Range("A27").Value = "C04"
$wsEmp.Range("B27").Value = "unemployment rate (%)"

$wsEmp.Range("A26").Value = "C03"
$wsEmp.Range("B26").Value = "employment/population ratio (%)"

$wsEmp.Range("A25").Value = "C02"
$wsEmp.Range("B25").Value = "labor force participation rate (%)"

$wsEmp.Range("A24").Value = "C01"
$wsEmp.Range("B24").Value = "total number"

# ---------------------------------------------------------------------------
# 2) Add the new "PubAssist" worksheet, placed after the last sheet ("Food").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPub = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPub.Name = "PubAssist"

# Header row
$wsPub.Range("A1").Value = "B19058"
$wsPub.Range("B1").Value = "Number"
$wsPub.Range("C1").Value = 2010
$wsPub.Range("D1").Value = 2011
$wsPub.Range("E1").Value = 2012
$wsPub.Range("F1").Value = 2013
$wsPub.Range("G1").Value = 2014
$wsPub.Range("H1").Value = 2015
$wsPub.Range("I1").Value = 2016
$wsPub.Range("J1").Value = 2017

# Row 2 / Row 3 labels first (A column), so the new shared strings are
# appended in the same order as the reference workbook: "total hh" then
# "hh with public assistance" before the reused "number of hh" label.
$wsPub.Range("A2").Value = "total hh"
$wsPub.Range("A3").Value = "hh with public assistance"

# Row 2: total hh
$wsPub.Range("B2").Value = 1
$wsPub.Range("C2").Value = "x"
$wsPub.Range("D2").Value = "x"
$wsPub.Range("E2").Value = "x"
$wsPub.Range("F2").Value = "x"
$wsPub.Range("G2").Value = "x"
$wsPub.Range("H2").Value = "x"
$wsPub.Range("I2").Value = "x"
$wsPub.Range("J2").Value = "x"
$wsPub.Range("K2").Value = "number of hh"

# Row 3: hh with public assistance
$wsPub.Range("B3").Value = 2
$wsPub.Range("C3").Value = "x"
$wsPub.Range("D3").Value = "x"
$wsPub.Range("E3").Value = "x"
$wsPub.Range("F3").Value = "x"
$wsPub.Range("G3").Value = "x"
$wsPub.Range("H3").Value = "x"
$wsPub.Range("I3").Value = "x"
$wsPub.Range("J3").Value = "x"
$wsPub.Range("K3").Value = "number of hh"

# Row 5: derived note
$wsPub.Range("K5").Value = "derive: % of hh with public assistance"

$wsPub.Range("K7").Select()

# ---------------------------------------------------------------------------
# 3) Restore per-sheet selections that moved in the edit.
# ---------------------------------------------------------------------------
$wsCensus = $wb.Worksheets.Item("Census")
$wsCensus.Activate()
$wsCensus.Range("D42").Select()

$wsEmp.Activate()
$wsEmp.Range("E16").Select()

$wsIncome = $wb.Worksheets.Item("Income")
$wsIncome.Activate()
$wsIncome.Range("B12").Select()

$wsEdu = $wb.Worksheets.Item("Edu")
$wsEdu.Activate()
$wsEdu.Range("Q19").Select()

$wsFood = $wb.Worksheets.Item("Food")
$wsFood.Activate()

# ---------------------------------------------------------------------------
# 4) PovFam becomes the active/selected sheet (activeTab moves from Edu to
#    PovFam), with its selection moved to I12.
# ---------------------------------------------------------------------------
$wsPovFam = $wb.Worksheets.Item("PovFam")
$wsPovFam.Activate()
$wsPovFam.Range("I12").Select()
